$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ndp -> Lgr4, Target cluster ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1127976666666667
$ws.Range("H2").Value = 0.338393
$ws.Range("M2").Value = 1.952294
$ws.Range("N2").Value = 5.856882
$ws.Range("O2").Value = 0.07575070565202183
$ws.Range("P2").Value = 0.07575070565202184
$ws.Range("Q2").Value = 0.2202142078473333
$ws.Range("R2").Value = 1.981927870626
$ws.Range("S2").Value = 0.07575070565202183
$ws.Range("T2").Value = 0.07575070565202184

# Row 3 (Ndp -> Lgr4, Target cluster FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1127976666666667
$ws.Range("H3").Value = 0.338393
$ws.Range("O3").Value = 0.5679402069281436
$ws.Range("P3").Value = 0.5679402069281437
$ws.Range("Q3").Value = 1.651053963086
$ws.Range("R3").Value = 14.859485667774
$ws.Range("S3").Value = 0.5679402069281436
$ws.Range("T3").Value = 0.5679402069281437

# Row 4 (Ndp -> Lgr4, Target cluster MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1127976666666667
$ws.Range("H4").Value = 0.338393
$ws.Range("M4").Value = 9.009963000000001
$ws.Range("N4").Value = 27.029889
$ws.Range("O4").Value = 0.3495944028658632
$ws.Range("P4").Value = 0.3495944028658634
$ws.Range("Q4").Value = 1.016302803153
$ws.Range("R4").Value = 9.146725228377001
$ws.Range("S4").Value = 0.3495944028658632
$ws.Range("T4").Value = 0.3495944028658634

# Row 5 (Ndp -> Lgr4, Target cluster Resolving-Mac)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1127976666666667
$ws.Range("H5").Value = 0.338393
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.173055
$ws.Range("N5").Value = 0.519165
$ws.Range("O5").Value = 0.006714684553971194
$ws.Range("P5").Value = 0.006714684553971196
$ws.Range("Q5").Value = 0.019520200205
$ws.Range("R5").Value = 0.175681801845
$ws.Range("S5").Value = 0.006714684553971194
$ws.Range("T5").Value = 0.006714684553971196
